$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header text tweak: "Done By:" -> "Done By"
$ws.Range("N2").Value = "Done By"

# 2) The "Identify 3 patterns..." review cards were miscategorised under
#    "Needs Reviewing" (column H:I) - move them over to "Done" (column J:K)
#    for rows 3-7, taking their value + fill colour + wrap/centre alignment
#    with them, and blank out the old H:I cells completely.
for ($r = 3; $r -le 7; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    $jCell = $ws.Cells.Item($r, 10)

    $jCell.Value = $hCell.Value2
    $jCell.Interior.Color = $hCell.Interior.Color
    $jCell.HorizontalAlignment = $hCell.HorizontalAlignment
    $jCell.WrapText = $hCell.WrapText

    $hCell.Value = $null
    $hCell.Interior.ColorIndex = -4142
    $hCell.Interior.Pattern = -4142
    $hCell.WrapText = $false
}

# 3) Tidy up the now-empty "Needs Reviewing" column: row 3 loses its border
#    entirely, rows 4-7 get a uniform thin box border around each cell.
$ws.Range("H3:I3").Borders.LineStyle = -4142
for ($r = 4; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Borders.LineStyle = 1
}

# 4) Widen the "Done" columns (J:K) to match the other content columns.
$ws.Range("J1:K1").EntireColumn.ColumnWidth = 16
